# Update the FHIR StructureDefinition spreadsheet:
#   - URL moved from ibm.com -> linuxforhealth.org (also updates the
#     Extension.url "Fixed Value" cell on the Elements sheet, which shares
#     the same text)
#   - Version bumped 7.0.0 -> 8.0.0
#   - Date bumped to the new publication timestamp
#   - Publisher renamed Alvearie Team -> LinuxForHealth Team
#   - The root Extension row's rolled-up Constraint(s) cell is cleared
#     (the ele-1/ext-1 constraint now only shows on the Extension.extension
#     child row)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: Property/Value table -------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/county-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: per-element grid ------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url row's "Fixed Value" column mirrors the StructureDefinition's
# own canonical URL - keep it in sync with the Metadata URL change above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/county-code"

# Root "Extension" row (row 2) no longer carries the rolled-up ele-1/ext-1
# constraint text in its Constraint(s) column (AI); it now only shows on
# the Extension.extension child row.
$elements.Range("AI2").Value = ""
